$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the last existing data row (row 70) down into the two
# new rows before writing values, so number formats / styles match.
$ws.Range("A70:I70").Copy() | Out-Null
$ws.Range("A71:I72").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Row 71
$ws.Range("A71").Value = 46037
$ws.Range("B71").Value = 5604
$ws.Range("C71").Value = 4305
$ws.Range("D71").Value = 3987
$ws.Range("E71").Value = 221
$ws.Range("F71").Value = 51
$ws.Range("G71").Value = 38
$ws.Range("H71").Value = 8
$ws.Range("I71").Value = 0

# Row 72
$ws.Range("A72").Value = 46038
$ws.Range("B72").Value = 5599
$ws.Range("C72").Value = 3981
$ws.Range("D72").Value = 3645
$ws.Range("E72").Value = 237
$ws.Range("F72").Value = 49
$ws.Range("G72").Value = 42
$ws.Range("H72").Value = 7
$ws.Range("I72").Value = 1

# Update selection / view to match the new active cell and scroll position.
$ws.Range("A72:I72").Select() | Out-Null
$ws.Application.ActiveWindow.ScrollRow = 55
